# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the last data row of the zh-cn and de-de report sheets,
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-18 08:24:57"
$wsZhCn.Range("G5").Value = "2016-02-18 08:25:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-18 08:25:08"
$wsDeDe.Range("G5").Value = "2016-02-18 08:26:20"
